$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows before the old row 12 (pushes old row12.. down to row15..) ---
$ws.Rows("11:13").Insert()

# --- New data block in rows 10-12, columns H-K (function generator wave table) ---
# Cell values are written in this specific order so the workbook's shared-string
# table is built up in the same sequence as the target file.
$ws.Range("I10").Value = "Sine wave mode"
$ws.Range("H11").Value = "Frequency"
$ws.Range("I11").Value = "0 - 65.789Hz; 526.32 Hz; 7KHz; 94.34KHz"
$ws.Range("J10").Value = "Triangular wave mode"
$ws.Range("J11").Value = "-"
$ws.Range("H12").Value = "Peak to Peak"
$ws.Range("J12").Value = "18V"
$ws.Range("K10").Value = "Square"
$ws.Range("I12").Value = "14V"
$ws.Range("K12").Value = "13.5V"
$ws.Range("K11").Value = "-"

$ws.Range("J11:K11").HorizontalAlignment = -4108

# --- New column widths for the added columns H, I, J ---
$ws.Columns("H").ColumnWidth = 14.449776785714286
$ws.Columns("I").ColumnWidth = 31.699776785714285
$ws.Columns("J").ColumnWidth = 17.867745535714285

# --- Update the view: selection moves to K13 ---
$ws.Range("K13").Select()
